$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.910.98'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.271.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.62%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.58'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.81'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.262.03'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.187'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -8.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.586'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.31'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -7.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.57'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '629.15'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.807.89'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.918.25'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.84'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.279.85'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.34'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.904'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '106.04'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.92'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -6.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.96'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -7.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.66'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.62'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.68'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.33'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.03'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.23'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.02'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.68%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '540.73'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.77%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.105'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.44'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.82%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.699.52'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.40'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0730'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -7.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.131'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.41'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '32.72'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.337'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -8.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.28'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0414'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.60'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -7.41%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.99%  '
